# Completed Ordering combo table across all devices and browsers.
# Fills in the device rows (Kindle Fire, iPad4, iPhone 6, Nexus 10, Galaxy S4)
# of the browser/device test-results matrix on Sheet1, columns C:U, and fixes
# up a few previously-entered values in row 13 (Galaxy S4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$colNums = @{
    "C" = 3;  "D" = 4;  "E" = 5;  "F" = 6;  "G" = 7;  "H" = 8;  "I" = 9;
    "J" = 10; "K" = 11; "L" = 12; "M" = 13; "N" = 14; "O" = 15; "P" = 16;
    "Q" = 17; "R" = 18; "S" = 19; "T" = 20; "U" = 21
}

$rowData = @{
    9 = @{
        "C" = "Yes"
        "D" = "No--price should be 7.49, not 7.94. Otherwise fine."
        "E" = "Fine, though way zoomed out."
        "F" = "The two can definitely be added, but the price should be 16.48, not 16.93."
        "G" = "Price wrong, but adds and removes pizzas fine. Zoomed out."
        "H" = "Price wrong, but adds and removes pizzas fine. Zoomed out."
        "I" = "Yes. Zoomed out."
        "J" = "No, same price problem. Zoomed out.Adds and removes fine."
        "K" = "Yes"
        "L" = "Price wrong; adds and removes fine."
        "M" = "Price wrong; adds and removes fine."
        "N" = "Price wrong; adds and removes fine."
        "O" = "Yes"
        "P" = "Wrong price; adds and removes fine."
        "Q" = "Wrong price; adds and removes fine."
        "R" = "Wrong price; adds and removes fine."
        "S" = "Wrong price; adds and removes fine."
        "T" = "Wrong price; adds and removes fine."
        "U" = "Wrong price; adds and removes fine."
    }
    10 = @{
        "C" = "Yes"
        "D" = "Price wrong. Otherwise fine."
        "E" = "Fine."
        "F" = "Price wrong. Fine, notwithstanding. (Maybe a juttery page?)"
        "G" = "Price wrong, but adds and removes pizzas fine."
        "H" = "Price wrong, but adds and removes pizzas fine."
        "I" = "Yes"
        "J" = "No, same price problem. Adds and removes fine."
        "K" = "Yes"
        "L" = "Price wrong; adds and removes fine."
        "M" = "Price wrong; adds and removes fine."
        "N" = "Price wrong; adds and removes fine."
        "O" = "Yes"
        "P" = "Wrong price; adds and removes fine."
        "Q" = "Wrong price; adds and removes fine."
        "R" = "Wrong price; adds and removes fine."
        "S" = "Wrong price; adds and removes fine."
        "T" = "Wrong price; adds and removes fine."
        "U" = "Wrong price; adds and removes fine."
    }
    11 = @{
        "C" = "Yes"
        "D" = "Price wrong. Images off plum. Otherwise fine."
        "E" = "Images offish. Otherwise fine."
        "F" = "Price wrong, and images off plum; needs optimised because jutters. Can still add."
        "G" = "Price wrong, but adds and removes pizzas fine. Zoomed in."
        "H" = "Price wrong, but adds and removes pizzas fine. Squashed."
        "I" = "Yes. Squashed."
        "J" = "No, same price problem. Zoomed in. Adds and removes fine."
        "K" = "Yes"
        "L" = "Price wrong; adds and removes fine."
        "M" = "Price wrong; adds and removes fine."
        "N" = "Price wrong; adds and removes fine."
        "O" = "Yes"
        "P" = "Wrong price; adds and removes fine."
        "Q" = "Wrong price; adds and removes fine."
        "R" = "Wrong price; adds and removes fine."
        "S" = "Wrong price; adds and removes fine."
        "T" = "Wrong price; adds and removes fine."
        "U" = "Wrong price; adds and removes fine."
    }
    12 = @{
        "C" = "Yes"
        "D" = "Price wrong. Otherwise fine."
        "E" = "Fine."
        "F" = "Price wrong but  can still add pizzas and remove. Fine, mostly."
        "G" = "Price wrong, but adds and removes pizzas fine."
        "H" = "Price wrong, but adds and removes pizzas fine. "
        "I" = "Yes."
        "J" = "No, same price problem. Adds and removes fine."
        "K" = "Yes"
        "L" = "Price wrong; adds and removes fine."
        "M" = "Price wrong; adds and removes fine."
        "N" = "Price wrong; adds and removes fine."
        "O" = "Yes"
        "P" = "Wrong price; adds and removes fine."
        "Q" = "Wrong price; adds and removes fine."
        "R" = "Wrong price; adds and removes fine."
        "S" = "Wrong price; adds and removes fine."
        "T" = "Wrong price; adds and removes fine."
        "U" = "Wrong price; adds and removes fine."
    }
    13 = @{
        "E" = "Fine."
        "G" = "No--should be 13.48, not 13.93. Images off plum. Fine for add"
        "H" = "Price wrong, but adds and removes pizzas fine. Squashed"
        "J" = "No--should be 25.47, not 25.92. Images off plum."
        "K" = "Yes"
        "L" = "Price wrong; adds and removes fine."
        "M" = "Price wrong; adds and removes fine."
        "N" = "Price wrong; adds and removes fine."
        "O" = "Yes"
        "P" = "Wrong price; adds and removes fine."
        "Q" = "Wrong price; adds and removes fine."
        "R" = "Wrong price; adds and removes fine."
        "S" = "Wrong price; adds and removes fine."
        "T" = "Wrong price; adds and removes fine."
        "U" = "Wrong price; adds and removes fine."
    }
}

foreach ($r in $rowData.Keys) {
    $cells = $rowData[$r]
    foreach ($colLetter in $cells.Keys) {
        $colNum = $colNums[$colLetter]
        $ws.Cells.Item($r, $colNum).Value = $cells[$colLetter]
    }
}

# Widen / autofit the columns now that the table is fully populated -- the
# author widened B by hand and best-fit the rest.
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).ColumnWidth = 71.33203125
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(5).AutoFit()
$ws.Columns.Item(6).AutoFit()
$ws.Columns.Item(7).AutoFit()
$ws.Columns.Item(8).AutoFit()
$ws.Columns.Item(9).AutoFit()
$ws.Columns.Item(10).AutoFit()
$ws.Columns.Item(11).AutoFit()
$ws.Columns.Item(12).AutoFit()
$ws.Range("M1:U1").Columns.AutoFit()

# Scroll down a row and select B24, matching the author's final cursor spot.
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$ws.Range("B24").Select()

# Restore the window geometry recorded for this save.
$win.Left = 3860
$win.Top = 1180
$win.Width = 22640
$win.Height = 15180
